{"js": "// Delete the table caption paragraph (\"Table S1. Summary of metagenome\n// data for Organic Lake samples.\") that precedes the table, leaving the\n// table as the first element of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst first = paragraphs.items[0];\nfirst.load(\"text\");\nawait context.sync();\n\nif (first.text.trim().indexOf(\"Table S1\") === 0) {\n  first.delete();\n  await context.sync();\n}\n", "ps1": "# Delete the table caption paragraph (\"Table S1. Summary of metagenome\n# data for Organic Lake samples.\") that precedes the table, leaving the\n# table as the first element of the document body.\n$d = $word.ActiveDocument\n\n$first = $d.Paragraphs(1)\nif ($first.Range.Text -like \"Table S1*\") {\n    $first.Range.Delete()\n}\n"}
